$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slang")

# Rename the existing shared string used by A8 from
# "ValidatePageNavigationOnVoiceInput" to "airtelAppEnglish"
$ws.Range("A8").Value = "airtelAppEnglish"

# Add a new "airtelAppHindi" test block (rows 12-13), mirroring the
# existing User/Password block pattern used by the other tests.
$ws.Range("A12").Value = "airtelAppHindi"

$ws.Range("B12").Value = "User"
$ws.Range("B12").Font.Bold = $true

$ws.Range("C12").Value = "Password"
$ws.Range("C12").Font.Bold = $true

$ws.Range("B13").Value = "sampleusername"
$ws.Range("B13").Style = "Hyperlink"

$ws.Range("C13").Value = "samplePwd"
$ws.Range("C13").Style = "Hyperlink"

$ws.Range("A13").Select()
